$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to store text (matches source workbook, where
# Price/Volume columns are inline strings, e.g. "1.00", "0.583", "64.612.22")
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '64.612.22'
$ws.Range("E2").Value = '  +2.60%  '

$ws.Range("D3").Value = '3.459.04'
$ws.Range("E3").Value = '  +2.72%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '576.28'
$ws.Range("E5").Value = '  +0.52%  '

$ws.Range("D6").Value = '158.30'
$ws.Range("E6").Value = '  +3.59%  '

$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").Value = '3.466.86'
$ws.Range("E8").Value = '  +2.84%  '

$ws.Range("D9").Value = '0.583'
$ws.Range("E9").Value = '  +11.11%  '

$ws.Range("E10").Value = '  -0.90%  '

$ws.Range("E11").Value = '  +5.36%  '

$ws.Range("D12").Value = '0.446'
$ws.Range("E12").Value = '  +2.50%  '

$ws.Range("D13").Value = '4.052.08'
$ws.Range("E13").Value = '  +2.71%  '

$ws.Range("E14").Value = '  -2.65%  '

$ws.Range("E15").Value = '  +7.65%  '

$ws.Range("D16").Value = '28.32'
$ws.Range("E16").Value = '  +5.23%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.538.14'
$ws.Range("E17").Value = '  +4.57%  '

$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '64.534.71'
$ws.Range("E18").Value = '  +2.38%  '

$ws.Range("E19").Value = '  +2.12%  '

$ws.Range("E20").Value = '  +4.00%  '

$ws.Range("D21").Value = '391.28'
$ws.Range("E21").Value = '  +1.71%  '

$ws.Range("D22").Value = '8.28'
$ws.Range("E22").Value = '  -1.13%  '

$ws.Range("D23").Value = '73.80'
$ws.Range("E23").Value = '  +5.00%  '

$ws.Range("D24").Value = '0.543'
$ws.Range("E24").Value = '  +1.93%  '

$ws.Range("D25").Value = '0.999'

$ws.Range("E26").Value = '  +24.91%  '

$ws.Range("D27").Value = '9.58'
$ws.Range("E27").Value = '  +3.46%  '

$ws.Range("E28").Value = '  +0.55%  '

$ws.Range("E29").Value = '  +0.25%  '

$ws.Range("D30").Value = '6.20'
$ws.Range("E30").Value = '  +11.15%  '

$ws.Range("E31").Value = '  +10.41%  '

$ws.Range("E32").Value = '  +0.67%  '

$ws.Range("D33").Value = '6.58'
$ws.Range("E33").Value = '  +3.65%  '

$ws.Range("D34").Value = '23.71'
$ws.Range("E34").Value = '  +3.03%  '

$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  +0.05%  '

$ws.Range("E36").Value = '  +5.07%  '

$ws.Range("E37").Value = '  +0.58%  '

$ws.Range("D38").Value = '161.11'
$ws.Range("E38").Value = '  +1.82%  '

$ws.Range("E39").Value = '  +1.67%  '

$ws.Range("E40").Value = '  +4.77%  '

$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").Value = '27.34'
$ws.Range("E41").Value = '  -0.08%  '

$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '2.937.54'
$ws.Range("E42").Value = '  +1.71%  '

$ws.Range("D43").Value = '0.0320'
$ws.Range("E43").Value = '  -2.16%  '

$ws.Range("D44").Value = '42.79'
$ws.Range("E44").Value = '  +4.89%  '

$ws.Range("D45").Value = '4.43'
$ws.Range("E45").Value = '  +4.83%  '

$ws.Range("D46").Value = '0.771'
$ws.Range("E46").Value = '  +3.03%  '

$ws.Range("D47").Value = '23.71'
$ws.Range("E47").Value = '  +8.40%  '

$ws.Range("E48").Value = '  +5.37%  '

$ws.Range("D49").Value = '2.26'
$ws.Range("E49").Value = '  +22.84%  '

$ws.Range("B50").Value = 'SuiNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D50").Value = '0.868'
$ws.Range("E50").Value = '  +7.79%  '

$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").Value = '0.108'
$ws.Range("E51").Value = '  +4.86%  '
